$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Solo en Casa" - add 5 new data rows (rows 2-6)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Solo en Casa")

$sheet1Rows = @(
    @(2, "Mesa y Lopez", "Z010100004", "CAMEL ACTIVA 10x20 Udes.", 8416500103257, 10, "2025-05-23 14:13:39", "left_only"),
    @(5, "Mesa y Lopez", "Z010100003", "CAMEL ACTIVA DOBLE 20x10 Udes.", 8416500103325, 6, "2025-05-23 14:19:20", "left_only"),
    @(1, "Mesa y Lopez", "Z010100003", "CAMEL ACTIVA DOBLE 20x10 Udes.", 8416500103325, 5, "2025-05-23 14:13:39", "left_only"),
    @(3, "Mesa y Lopez", "Z010100020", "CAMEL RYO LEGEND 10X30 Grms.", 8416500021797, 15, "2025-05-23 14:13:39", "left_only"),
    @(4, "Mesa y Lopez", "Z010200000", "WINSTON CLAS.BOX CT 10x20 Ude.", 8416500140948, 20, "2025-05-23 14:13:39", "left_only")
)

$r = 2
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $ws1.Cells.Item($r, 11).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet 2: "Solo en AS400" - insert new row at position 9
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Solo en AS400")

$ws2.Rows.Item(9).Insert()
$ws2.Cells.Item(9, 4).Value = "CAMEL BLUE CARTON 10x20 Ud."
$ws2.Cells.Item(9, 8).Value = 2011000001
$ws2.Cells.Item(9, 9).Value = "720 12103"
$ws2.Cells.Item(9, 10).Value = 16
$ws2.Cells.Item(9, 11).Value = "right_only"

# ---------------------------------------------------------------------------
# Sheet 3: "Diferencias Stock" - replace row 2 content, remove row 3
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Diferencias Stock")

$ws3.Cells.Item(2, 1).Value = 6
$ws3.Cells.Item(2, 2).Value = "Mesa y Lopez"
$ws3.Cells.Item(2, 3).Value = "Z010100000"
$ws3.Cells.Item(2, 4).Value = "CAMEL FILTERS CARTON 10x20 Ud."
$ws3.Cells.Item(2, 5).Value = 8416500140962
$ws3.Cells.Item(2, 6).Value = 7
$ws3.Cells.Item(2, 7).Value = "2025-05-23 14:22:25"
$ws3.Cells.Item(2, 8).Value = 2011000000
$ws3.Cells.Item(2, 9).Value = "720 12100"
$ws3.Cells.Item(2, 10).Value = 17
$ws3.Cells.Item(2, 11).Value = "both"
$ws3.Cells.Item(2, 12).Value = -10

$ws3.Rows.Item(3).Delete()
